$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: G7 already has style s="4" (volunteer-style empty cell) -- just fill the value.
$ws.Range("G7").Value = 1.0

# Row 7: H7 is a brand-new cell that should look like the existing "/" sign-off cells
# (e.g. F4), which carry style s="3" and the shared string "/". Copy F4's formatting
# so the same cellXf gets reused instead of generating a new one, then set the value.
$ws.Range("F4").Copy()
$ws.Range("H7").PasteSpecial(-4122)
$ws.Range("H7").Value = "/"

# Row 10: G10 is a new cell that should match F10's style (s="3") with value 1.0.
$ws.Range("F10").Copy()
$ws.Range("G10").PasteSpecial(-4122)
$ws.Range("G10").Value = 1.0

# Row 10: H10 is a new "/" sign-off cell, same pattern as H7/F4.
$ws.Range("F4").Copy()
$ws.Range("H10").PasteSpecial(-4122)
$ws.Range("H10").Value = "/"

$excel.CutCopyMode = 0
